$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds mixed numeric-looking text (e.g. "1.002", "251.06") that must
# stay as literal text, matching the source data which stores everything as
# inline strings. Pre-format the whole data range as Text so Excel does not
# auto-convert the new values into numbers when we assign them below.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "30.338.25"
$ws.Range("E2").Value = "  +0.53%  "
$ws.Range("D3").Value = "1.932.60"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "251.06"
$ws.Range("E5").Value = "  +1.97%  "
$ws.Range("D6").Value = "0.7174"
$ws.Range("E6").Value = "  +0.20%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "0.3274"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("D9").Value = "27.50"
$ws.Range("E9").Value = "  +4.77%  "
$ws.Range("D10").Value = "0.07178"
$ws.Range("E10").Value = "  +5.06%  "
$ws.Range("D11").Value = "0.8007"
$ws.Range("E11").Value = "  +0.93%  "
$ws.Range("D12").Value = "0.08078"
$ws.Range("E12").Value = "  +2.00%  "
$ws.Range("D13").Value = "1.932.37"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("D14").Value = "5.414"
$ws.Range("E14").Value = "  +0.55%  "
$ws.Range("D15").Value = "94.43"
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").Value = "14.84"
$ws.Range("E16").Value = "  +2.72%  "
$ws.Range("D17").Value = "30.321.26"
$ws.Range("D18").Value = "252.00"
$ws.Range("E18").Value = "  -2.71%  "
$ws.Range("D19").Value = "0.000008118"
$ws.Range("E19").Value = "  +2.68%  "
$ws.Range("D20").Value = "5.794"
$ws.Range("E20").Value = "  -0.61%  "
$ws.Range("D21").Value = "2.186.01"
$ws.Range("E21").Value = "  +0.70%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("D24").Value = "6.919"
$ws.Range("E24").Value = "  +1.11%  "
$ws.Range("D25").Value = "9.715"
$ws.Range("E25").Value = "  +0.32%  "
$ws.Range("D26").Value = "165.43"
$ws.Range("D27").Value = "19.23"
$ws.Range("E27").Value = "  +2.30%  "
$ws.Range("D28").Value = "2.329"
$ws.Range("E28").Value = "  +4.17%  "
$ws.Range("D29").Value = "0.1290"
$ws.Range("E29").Value = "  -2.93%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").Value = "1.542"
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("D32").Value = "4.417"
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").Value = "4.205"
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("D34").Value = "0.05196"
$ws.Range("E34").Value = "  +3.20%  "
$ws.Range("D35").Value = "1.265"
$ws.Range("E35").Value = "  +6.36%  "
$ws.Range("D36").Value = "0.7466"
$ws.Range("E36").Value = "  +1.20%  "
$ws.Range("D37").Value = "2.768"
$ws.Range("E37").Value = "  +1.32%  "
$ws.Range("D38").Value = "0.01961"
$ws.Range("E38").Value = "  +0.96%  "
$ws.Range("D39").Value = "2.798"
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("D40").Value = "78.86"
$ws.Range("E40").Value = "  -1.23%  "
$ws.Range("D41").Value = "6.455"
$ws.Range("E41").Value = "  -0.41%  "
$ws.Range("D42").Value = "0.4526"
$ws.Range("E42").Value = "  +2.53%  "
$ws.Range("D43").Value = "2.023"
$ws.Range("E43").Value = "  +1.12%  "
$ws.Range("E44").Value = "  +0.17%  "
$ws.Range("D45").Value = "0.8408"
$ws.Range("E45").Value = "  +1.41%  "
$ws.Range("D46").Value = "101.81"
$ws.Range("E46").Value = "  -0.49%  "
$ws.Range("D47").Value = "9.795"
$ws.Range("E47").Value = "  +1.23%  "
$ws.Range("D48").Value = "7.401"
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("D49").Value = "36.61"
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("D50").Value = "0.4177"
$ws.Range("E50").Value = "  +2.03%  "
$ws.Range("D51").Value = "0.06063"
$ws.Range("E51").Value = "  +2.40%  "
